# Generate Report for Handoff
# The localization of 819a2cbc-9f1b-4372-9acc-939d459cfb78.md has been
# re-handed-off: update its status from "Handed back: in sync with en-US"
# to "Ready for handoff" on all three sheets, and record the new handoff
# timestamps on the per-locale detail sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = "Ready for handoff"
$zhcn.Range("D3").Value = "2016-03-02 10:04:33"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = "Ready for handoff"
$dede.Range("D3").Value = "2016-03-02 10:04:44"
